{"js": "// Move the \"_GoBack\" bookmark from its old location (between the \"P\" run and\n// the \"lease send...\" run, inside the closing paragraph) to a new location\n// immediately after the \"Patient.City\" merge-field text in the address block.\n\n// 1) Find the \"Patient.City\" text in the document body and collapse a range\n//    to its end so we can insert the bookmark right after it.\nconst cityResults = context.document.body.search(\"Patient.City\", { matchCase: true, matchWholeWord: false });\ncityResults.load(\"items\");\nawait context.sync();\n\nif (cityResults.items.length === 0) {\n  throw new Error(\"Could not find 'Patient.City' text in the document.\");\n}\n\nconst cityRange = cityResults.items[0];\nconst afterCity = cityRange.getRange(\"End\");\n\n// 2) Remove the existing \"_GoBack\" bookmark (if present) before re-inserting\n//    it at the new location, since a bookmark name must be unique.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 3) Insert the bookmark right after \"Patient.City\".\nafterCity.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Move the \"_GoBack\" bookmark from its old location (between the \"P\" run and\n# the \"lease send...\" run, inside the closing paragraph) to a new location\n# immediately after the \"Patient.City\" merge-field text in the address block.\n\n$d = $word.ActiveDocument\n\n# 1) Find the \"Patient.City\" text in the document body.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"Patient.City\")\nif (-not $found) {\n    throw \"Could not find 'Patient.City' text in the document.\"\n}\n\n# 2) Collapse the found range to its end (a zero-length caret right after\n#    \"Patient.City\").\n$rng.Collapse(0)  # wdCollapseEnd\n\n# 3) Re-adding a bookmark with an existing name moves it to the new range,\n#    so this both removes the old \"_GoBack\" bookmark and creates it here.\n$d.Bookmarks.Add(\"_GoBack\", $rng)\n"}
